$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.064.18'
$ws.Range('E2').Value = '  +0.65%  '
$ws.Range('D3').Value = '2.377.10'
$ws.Range('E3').Value = '  +2.42%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '301.98'
$ws.Range('E5').Value = '  -0.09%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '96.86'
$ws.Range('E7').Value = '  -0.55%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('E9').Value = '  +1.41%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.21'
$ws.Range('E10').Value = '  -0.95%  '
$ws.Range('E11').Value = '  +0.43%  '
$ws.Range('E12').Value = '  +2.72%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.31'
$ws.Range('E13').Value = '  -4.75%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.80'
$ws.Range('E14').Value = '  +0.81%  '
$ws.Range('D15').Value = '2.750.23'
$ws.Range('D16').Value = '2.395.44'
$ws.Range('E16').Value = '  +4.02%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.804'
$ws.Range('E17').Value = '  +2.26%  '
$ws.Range('D18').Value = '43.030.67'
$ws.Range('E18').Value = '  +0.72%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.18'
$ws.Range('E19').Value = '  -0.59%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.32'
$ws.Range('E21').Value = '  -0.36%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.14'
$ws.Range('E22').Value = '  +0.18%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.55'
$ws.Range('E23').Value = '  -0.21%  '
$ws.Range('E24').Value = '  -1.57%  '
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.44'
$ws.Range('E25').Value = '  +1.03%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.07%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.97'
$ws.Range('E27').Value = '  +2.36%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.36'
$ws.Range('E28').Value = '  -0.20%  '
$ws.Range('E29').Value = '  +1.27%  '
$ws.Range('E30').Value = '  -3.06%  '
$ws.Range('E31').Value = '  -0.01%  '
$ws.Range('E32').Value = '  +1.16%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0742'
$ws.Range('E33').Value = '  +5.71%  '
$ws.Range('E34').Value = '  -1.39%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.86'
$ws.Range('E35').Value = '  +4.82%  '
$ws.Range('E36').Value = '  +5.15%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.34'
$ws.Range('E37').Value = '  -2.46%  '
$ws.Range('E38').Value = '  -0.69%  '
$ws.Range('E39').Value = '  +14.24%  '
$ws.Range('E40').Value = '  +4.00%  '
$ws.Range('E41').Value = '  -0.58%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '113.86'
$ws.Range('E42').Value = '  -31.56%  '
$ws.Range('D43').Value = '1.952.56'
$ws.Range('E43').Value = '  +0.06%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0280'
$ws.Range('E44').Value = '  +0.52%  '
$ws.Range('E45').Value = '  +2.34%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.74'
$ws.Range('E46').Value = '  -0.14%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.18'
$ws.Range('E47').Value = '  -11.69%  '
$ws.Range('D48').Value = '2.613.03'
$ws.Range('E48').Value = '  +2.63%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.52'
$ws.Range('E49').Value = '  +2.86%  '
$ws.Range('B50').Value = 'MultiversX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '52.40'
$ws.Range('E50').Value = '  -1.95%  '
$ws.Range('B51').Value = 'BitcoinSV'
$ws.Range('C51').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.52'
$ws.Range('E51').Value = '  +0.62%  '
